$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G1").Value = "Test"
